# "Generate Report for Handoff"
#
# A new handoff job (9c8dc578-c8fd-4a86-8163-eda76562a86e, content hash
# d018b6858bd87d1d6e28e1db06812cf978e822ba) replaces the old one
# (51b7e41d-8619-4e7a-954f-cc1905cf15c4 / 7a08454758c14b493fb3f597db47f948e52c7d95),
# and the stale "e8115f20-2b93-41f6-9a74-8fe89d9ff1bd.md / Handoff transform
# failed" row is dropped from every sheet (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$oldUuid = "51b7e41d-8619-4e7a-954f-cc1905cf15c4"
$newUuid = "9c8dc578-c8fd-4a86-8163-eda76562a86e"
$oldHash = "7a08454758c14b493fb3f597db47f948e52c7d95"
$newHash = "d018b6858bd87d1d6e28e1db06812cf978e822ba"

# ---------------------------------------------------------------------
# Sheet "Overview": drop row 3 (the removed file), rewrite the md UUID.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "$newUuid.md"
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

$ws1.Range("A3").Value = "e8115f20-2b93-41f6-9a74-8fe89d9ff1bd.md"
$ws1.Range("B3").Value = "Handoff transform failed"
$ws1.Range("C3").Value = "Handoff transform failed"

$ws1.Range("A4").Value = ".localization-config"
$ws1.Range("B4").Value = "Not to be localized"
$ws1.Range("C4").Value = "Not to be localized"

$ws1.Rows.Item(3).Delete()

$ws1.Hyperlinks.Delete()
$h = $ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cf28a08d29dc27d2cba740277ce6969429d2de91/e2e/$newUuid.md", "", "", "$newUuid.md")
$ws1.Range("A2").Style = "HyperLink"
$h = $ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/cf28a08d29dc27d2cba740277ce6969429d2de91/.localization-config", "", "", ".localization-config")
$ws1.Range("A3").Style = "HyperLink"

# ---------------------------------------------------------------------
# Sheet "zh-cn": drop row 3, rewrite UUID/hash/timestamp for row 2.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "$newUuid.md"
$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("C2").Value = "$newUuid.$newHash.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-02-18 08:12:39"
$ws2.Range("G2").Value = "0001-01-01 00:00:00"
$ws2.Range("H2").Value = "Include"

$ws2.Range("A3").Value = "e8115f20-2b93-41f6-9a74-8fe89d9ff1bd.md"
$ws2.Range("B3").Value = "Handoff transform failed"
$ws2.Range("D3").Value = "0001-01-01 00:00:00"
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Ignored"

$ws2.Range("A4").Value = ".localization-config"
$ws2.Range("B4").Value = "Not to be localized"
$ws2.Range("D4").Value = "0001-01-01 00:00:00"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "Ignored"

$ws2.Rows.Item(3).Delete()

$ws2.Hyperlinks.Delete()
$h = $ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cf28a08d29dc27d2cba740277ce6969429d2de91/e2e/$newUuid.md", "", "", "$newUuid.md")
$ws2.Range("A2").Style = "HyperLink"
$h = $ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6365b8b2ec7b9b69a6d39aca0834138dcbb41276/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$newUuid.$newHash.zh-cn.xlf", "", "", "$newUuid.$newHash.zh-cn.xlf")
$ws2.Range("C2").Style = "HyperLink"
$h = $ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/cf28a08d29dc27d2cba740277ce6969429d2de91/.localization-config", "", "", ".localization-config")
$ws2.Range("A3").Style = "HyperLink"

# ---------------------------------------------------------------------
# Sheet "de-de": drop row 3, rewrite UUID/hash/timestamp for row 2.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "$newUuid.md"
$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("C2").Value = "$newUuid.$newHash.de-de.xlf"
$ws3.Range("D2").Value = "2016-02-18 08:12:49"
$ws3.Range("G2").Value = "0001-01-01 00:00:00"
$ws3.Range("H2").Value = "Include"

$ws3.Range("A3").Value = "e8115f20-2b93-41f6-9a74-8fe89d9ff1bd.md"
$ws3.Range("B3").Value = "Handoff transform failed"
$ws3.Range("D3").Value = "0001-01-01 00:00:00"
$ws3.Range("G3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Ignored"

$ws3.Range("A4").Value = ".localization-config"
$ws3.Range("B4").Value = "Not to be localized"
$ws3.Range("D4").Value = "0001-01-01 00:00:00"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "Ignored"

$ws3.Rows.Item(3).Delete()

$ws3.Hyperlinks.Delete()
$h = $ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cf28a08d29dc27d2cba740277ce6969429d2de91/e2e/$newUuid.md", "", "", "$newUuid.md")
$ws3.Range("A2").Style = "HyperLink"
$h = $ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/666137a6858019e0f5156253e6b1dc0cda57fd87/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$newUuid.$newHash.de-de.xlf", "", "", "$newUuid.$newHash.de-de.xlf")
$ws3.Range("C2").Style = "HyperLink"
$h = $ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/cf28a08d29dc27d2cba740277ce6969429d2de91/.localization-config", "", "", ".localization-config")
$ws3.Range("A3").Style = "HyperLink"
